$d = $word.ActiveDocument

# Move to the very end of the document content (after the last paragraph,
# which currently holds the Bibliografia text) and open up a fresh blank
# paragraph there so InsertXML has a target paragraph to replace/expand.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter() | Out-Null

$targetPara = $d.Paragraphs.Last
$targetRange = $targetPara.Range

$insertXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Requisitos</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>LOB1055 -  Fundamentos de Engenharia de Segurança no Trabalho  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4076 -  Termodinâmica Aplicada  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4095 -  Química Geral Experimental  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4245 -  Ergonomia  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4251 -  Fundamentos de Química  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4252 -  Fundamentos de Fenômenos de Transporte  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4260 -  Controle Estatístico da Qualidade  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4262 -  Automação e Controle  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4263 -  Planejamento e Gestão da Manutenção  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4270 -  Planejamento, Programação e Controle da Produção II  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4271 -  Planejamento de Experimentos  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4272 -  Projeto da Fabrica  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1006 -  Cálculo IV  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1040 -  Laboratório de Eletricidade  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1046 -  Engenharia do Meio Ambiente  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1053 -  Física III  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4213 -  Contabilidade e Custos  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4236 -  Projeto Integrado de Engenharia de Produção I  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4237 -  Projeto Integrado de Engenharia de Produção II  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4238 -  Projeto Integrado em Engenharia de Produção III  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4241 -  Sistemas de Apoio à Decisão  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4255 -  Inovação Tecnológica  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4257 -  Gestão de Projetos  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4258 -  Pesquisa Operacional II  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4259 -  Processos de Desenvolvimento de Serviços  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4261 -  Planejamento, Programação e Controle da Produção I  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1003 -  Cálculo I  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1011 -  Eletricidade Aplicada  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1012 -  Estatística  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1024 -  Mecânica  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1036 -  Geometria Analítica  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1038 -  Física Experimental I  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1039 -  Física Experimental III  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1041 -  Física Experimental II  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1052 -  Cálculo III  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4201 -  Introdução à Engenharia de Produção  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4203 -  Sistemas Produtivos  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4239 -  Administração e Organização I  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4240 -  Administração e Organização II  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4253 -  Processos Químicos Industriais  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4264 -  Engenharia da Sustentabilidade  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1004 -  Cálculo II  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1018 -  Física I  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1019 -  Física II  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1049 -  Estatística Multivariada  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4204 -  Economia Geral  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4206 -  Pesquisa Operacional I  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4209 -  Engenharia da Qualidade  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4222 -  Engenharia Econômica e Finanças  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4234 -  Empreendedorismo e Inovação  (Requisito)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$targetRange.InsertXML($insertXml) | Out-Null

# InsertXML leaves one extra blank paragraph mark behind (the original
# paragraph mark that used to terminate the document). Merge it away so
# the new "ListBullet" paragraph becomes the final paragraph again, same
# as before the edit.
$trailingEmpty = $d.Paragraphs.Last
$trailingIndex = $trailingEmpty.Index
$prevPara = $d.Paragraphs.Item($trailingIndex - 1)

$cleanupRange = $d.Range($prevPara.Range.End - 1, $trailingEmpty.Range.End)
$cleanupRange.Delete()

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
Write-Output ("Last paragraph style=" + $d.Paragraphs.Last.Style)
Write-Output ("Last paragraph text=" + $d.Paragraphs.Last.Range.Text)
